# Add "Power Level" and "Class" variables (columns C and D) to the raw data sheet.
# Column C ("Power Level") is filled in first, top to bottom, then column D ("Class"),
# mirroring the order in which the new values were authored (and thus the order new
# strings land in the shared string table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: Power Level ---
$ws.Range("C1").Value = "Power Level"
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 11000
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 4.2
$ws.Range("C7").Value = 150
$ws.Range("C8").Value = 150
$ws.Range("C9").Value = 210
$ws.Range("C10").Value = 1450
$ws.Range("C11").Value = 710
$ws.Range("C12").Value = 6
# Text value starting with ">" must stay text, not be mis-parsed
$ws.Range("C13").Value = "> 9000"
$ws.Range("C14").Value = 12
$ws.Range("C15").Value = 0.01

# --- Column D: Class ---
$ws.Range("D1").Value = "Class"
$ws.Range("D2").Value = "Wizard"
$ws.Range("D3").Value = "Wizard"
$ws.Range("D4").Value = "Thief"
$ws.Range("D5").Value = "Warrior"
$ws.Range("D6").Value = "Thief"
$ws.Range("D7").Value = "Wizard"
$ws.Range("D8").Value = "Wizard"
$ws.Range("D9").Value = "Wizard"
$ws.Range("D10").Value = "Fighter"
$ws.Range("D11").Value = "Fighter"
$ws.Range("D12").Value = "Thief"
$ws.Range("D13").Value = "Bard"
$ws.Range("D14").Value = "Wizard"
$ws.Range("D15").Value = "Fighter"

# Update the selection to match the final state of the workbook
$ws.Range("E8").Select()
